# Apply updated symbol-list data (prices, 1h volume %, and rotated
# coin listing for rows 6-19) as scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="323.44"; E="-1.93%" }
    @{ Row=3; E="-1.19%" }
    @{ Row=4; D="5.879"; E="11.62%" }
    @{ Row=5; D="0.08032"; E="-0.80%" }
    @{ Row=6; B="KuCoinToken"; C="https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"; D="8.650"; E="-0.03%" }
    @{ Row=7; B="FTXToken"; C="https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; D="1.947"; E="0.93%" }
    @{ Row=8; B="BTSEToken"; C="https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; D="2.950"; E="-0.25%" }
    @{ Row=9; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="0.9288"; E="-0.78%" }
    @{ Row=10; B="LiechtensteinCryptoassetsExchange"; C="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; D="0.1279"; E="-5.03%" }
    @{ Row=11; B="WazirX"; C="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D="0.1960"; E="-0.73%" }
    @{ Row=12; B="MCDex"; C="https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; D="8.722"; E="34.35%" }
    @{ Row=13; B="MandalaExchangeToken"; C="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D="0.09195"; E="1.06%" }
    @{ Row=14; B="BitrueCoin"; C="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D="0.03555"; E="1.48%" }
    @{ Row=15; B="BitMartToken"; C="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D="0.1046"; E="9.16%" }
    @{ Row=16; B="BitForexToken"; C="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D="0.001296"; E="-7.98%" }
    @{ Row=17; B="TigerCash"; C="https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D="0.006088"; E="-1.14%" }
    @{ Row=18; B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="3.349"; E="-1.00%" }
    @{ Row=19; B="GateToken"; C="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; D="4.575"; E="1.18%" }
    @{ Row=20; D="0.3537"; E="0.52%" }
    @{ Row=21; E="4.39%" }
    @{ Row=23; D="0.04403"; E="-1.01%" }
    @{ Row=24; D="0.001263"; E="3.29%" }
    @{ Row=25; D="0.004397"; E="1.78%" }
    @{ Row=26; E="-11.72%" }
    @{ Row=39; D="0.02538"; E="1.27%" }
    @{ Row=40; D="0.05267"; E="1.44%" }
    @{ Row=41; D="0.007451"; E="-3.39%" }
    @{ Row=42; D="0.009615"; E="4.57%" }
    @{ Row=43; E="-1.56%" }
    @{ Row=44; D="0.002117"; E="-2.12%" }
    @{ Row=45; D="0.009992"; E="10.90%" }
    @{ Row=46; D="0.00006720"; E="1.38%" }
    @{ Row=47; E="-0.06%" }
    @{ Row=48; D="0.003003"; E="-10.25%" }
    @{ Row=49; E="-7.73%" }
    @{ Row=50; E="-0.06%" }
    @{ Row=51; E="-0.06%" }
)

foreach ($u in $updates) {
    foreach ($col in @("B", "C", "D", "E")) {
        if ($u.ContainsKey($col)) {
            $ref = "$col$($u.Row)"
            if ($col -eq "D" -or $col -eq "E") {
                # Price / Volume(1h) columns hold literal text (trailing
                # zeros, "--", "%") that must survive verbatim -- force
                # text format so Excel does not coerce them to numbers.
                $ws.Range($ref).NumberFormat = "@"
            }
            $ws.Range($ref).Value = $u[$col]
        }
    }
}
